$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2192307692307692
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.01923076923076923
$ws.Range("P2").Value = 0.1346153846153846
$ws.Range("S2").Value = 0.1269230769230769
$ws.Range("C3").Value = 0.01515151515151515
$ws.Range("J3").Value = 0.05303030303030303
$ws.Range("P3").Value = 0.6742424242424242
$ws.Range("S3").Value = 0.2575757575757576
$ws.Range("J4").Value = 0.05405405405405406
$ws.Range("P4").Value = 0.5135135135135135
$ws.Range("S4").Value = 0.4324324324324325
$ws.Range("B6").Value = 0.08163265306122448
$ws.Range("D6").Value = 0.00510204081632653
$ws.Range("F6").Value = 0.06122448979591837
$ws.Range("J6").Value = 0.2551020408163265
$ws.Range("O6").Value = 0.01530612244897959
$ws.Range("Q6").Value = 0.1275510204081633
$ws.Range("R6").Value = 0.07653061224489796
$ws.Range("S6").Value = 0.3775510204081632
$ws.Range("B7").Value = 0.08823529411764706
$ws.Range("D7").Value = 0.02352941176470588
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1529411764705882
$ws.Range("O7").Value = 0.02352941176470588
$ws.Range("Q7").Value = 0.1235294117647059
$ws.Range("R7").Value = 0.08235294117647059
$ws.Range("S7").Value = 0.4470588235294118
$ws.Range("B8").Value = 0.09560723514211886
$ws.Range("D8").Value = 0.01550387596899225
$ws.Range("F8").Value = 0.05426356589147287
$ws.Range("J8").Value = 0.124031007751938
$ws.Range("O8").Value = 0.020671834625323
$ws.Range("Q8").Value = 0.1808785529715762
$ws.Range("R8").Value = 0.1059431524547804
$ws.Range("S8").Value = 0.4031007751937984
$ws.Range("B9").Value = 0.07051282051282051
$ws.Range("D9").Value = 0.01923076923076923
$ws.Range("F9").Value = 0.07692307692307693
$ws.Range("J9").Value = 0.1858974358974359
$ws.Range("O9").Value = 0.01282051282051282
$ws.Range("Q9").Value = 0.1730769230769231
$ws.Range("S9").Value = 0.3782051282051282
$ws.Range("B10").Value = 0.1099476439790576
$ws.Range("D10").Value = 0.02006980802792321
$ws.Range("E10").Value = 0.0008726003490401396
$ws.Range("F10").Value = 0.07504363001745201
$ws.Range("J10").Value = 0.1273996509598604
$ws.Range("O10").Value = 0.01483420593368237
$ws.Range("Q10").Value = 0.1849912739965096
$ws.Range("R10").Value = 0.08813263525305411
$ws.Range("S10").Value = 0.3787085514834206
$ws.Range("G11").Value = 0.1691729323308271
$ws.Range("J11").Value = 0.08646616541353383
$ws.Range("K11").Value = 0.2105263157894737
$ws.Range("L11").Value = 0.5263157894736842
$ws.Range("S11").Value = 0.007518796992481203
$ws.Range("G12").Value = 0.6887417218543046
$ws.Range("J12").Value = 0.1920529801324503
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.0728476821192053
$ws.Range("S12").Value = 0.03973509933774835
$ws.Range("G13").Value = 0.6739130434782609
$ws.Range("J13").Value = 0.3260869565217391
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.6666666666666666
$ws.Range("F15").Value = 0.04419889502762431
$ws.Range("H15").Value = 0.1657458563535912
$ws.Range("I15").Value = 0.03867403314917127
$ws.Range("J15").Value = 0.3591160220994475
$ws.Range("K15").Value = 0.07734806629834254
$ws.Range("O15").Value = 0.04419889502762431
$ws.Range("S15").Value = 0.2707182320441989
$ws.Range("F16").Value = 0.007092198581560284
$ws.Range("H16").Value = 0.1631205673758865
$ws.Range("J16").Value = 0.4893617021276596
$ws.Range("K16").Value = 0.0851063829787234
$ws.Range("M16").Value = 0.03546099290780142
$ws.Range("N16").Value = 0.007092198581560284
$ws.Range("S16").Value = 0.1276595744680851
$ws.Range("F17").Value = 0.01420454545454545
$ws.Range("H17").Value = 0.1903409090909091
$ws.Range("I17").Value = 0.09943181818181818
$ws.Range("J17").Value = 0.4261363636363636
$ws.Range("K17").Value = 0.1051136363636364
$ws.Range("M17").Value = 0.03125
$ws.Range("N17").Value = 0.002840909090909091
$ws.Range("O17").Value = 0.04545454545454546
$ws.Range("S17").Value = 0.08522727272727272
$ws.Range("F18").Value = 0.01657458563535912
$ws.Range("H18").Value = 0.2596685082872928
$ws.Range("I18").Value = 0.06629834254143646
$ws.Range("J18").Value = 0.3867403314917127
$ws.Range("K18").Value = 0.1049723756906077
$ws.Range("M18").Value = 0.01657458563535912
$ws.Range("O18").Value = 0.04972375690607735
$ws.Range("S18").Value = 0.09944751381215469
$ws.Range("F19").Value = 0.01888489208633094
$ws.Range("H19").Value = 0.2014388489208633
$ws.Range("I19").Value = 0.08723021582733813
$ws.Range("J19").Value = 0.3803956834532374
$ws.Range("K19").Value = 0.1106115107913669
$ws.Range("M19").Value = 0.02338129496402878
$ws.Range("N19").Value = 0.0008992805755395684
$ws.Range("O19").Value = 0.07194244604316546
$ws.Range("S19").Value = 0.1052158273381295
